$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "Programa resumido:" ... "Requisitos:" block (rows 12-22)
# down by 3 rows (to rows 15-25) to make room for the new "Docentes
# responsáveis:" block. Values are copied directly (bottom-up, so a row is
# always fully read before anything overwrites it) instead of using
# Rows.Insert(), because Insert() also stamps neighbouring, otherwise-empty
# columns of the shifted rows with inherited formatting that shouldn't be
# there.
for ($oldRow = 22; $oldRow -ge 12; $oldRow--) {
    $newRow = $oldRow + 3

    # .Text is read eagerly (unlike .Value, which can resolve lazily against
    # the live cell and pick up later writes), so it's safe even though
    # source and destination ranges overlap.
    $aText = $ws.Cells.Item($oldRow, 1).Text
    $bText = $ws.Cells.Item($oldRow, 2).Text
    $cText = $ws.Cells.Item($oldRow, 3).Text

    if ([string]::IsNullOrEmpty($aText)) {
        $ws.Cells.Item($newRow, 1).ClearContents()
    } else {
        $ws.Cells.Item($newRow, 1).Value = $aText
    }

    if ([string]::IsNullOrEmpty($bText)) {
        $ws.Cells.Item($newRow, 2).ClearContents()
    } else {
        $ws.Cells.Item($newRow, 2).Value = $bText
    }

    if ([string]::IsNullOrEmpty($cText)) {
        $ws.Cells.Item($newRow, 3).ClearContents()
    } else {
        $ws.Cells.Item($newRow, 3).Value = $cText
    }

    $ws.Rows.Item($newRow).RowHeight = $ws.Rows.Item($oldRow).RowHeight
}

# Rows 12-14 are reused for the new "Docentes responsáveis:" block, so clear
# any of their old content that wouldn't otherwise be overwritten below.
$ws.Range("A12:C14").ClearContents()

# Populate the new rows 12-14 with the "Docentes responsáveis:" block.
$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Range("B13").Value = "3480026 - João Paulo Pascon"
$ws.Range("C13").Value = "3480026 - João Paulo Pascon"
$ws.Range("B14").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C14").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
